# Update raw loudness data (column G) with corrected values.
# Column I (DiffThreshold = G - H) recalculates automatically via its formula.
$wb = $excel.ActiveWorkbook
$ws = $wb.ActiveSheet

$vals = New-Object 'object[,]' 8,1
$vals[0,0] = 16.195573146292585
$vals[1,0] = 15.637605210420842
$vals[2,0] = 16.430507014028059
$vals[3,0] = 16.026714428857712
$vals[4,0] = 15.211787575150302
$vals[5,0] = 16.775566132264526
$vals[6,0] = 18.449469939879766
$vals[7,0] = 16.540632264529059
$ws.Range("G2:G9").Value = $vals

$vals = New-Object 'object[,]' 63,1
$vals[0,0] = 15.982664328657314
$vals[1,0] = 20.196790581162325
$vals[2,0] = 15.931272545090176
$vals[3,0] = 17.49505110220441
$vals[4,0] = 17.172017034068141
$vals[5,0] = 17.127966933867736
$vals[6,0] = 18.677062124248497
$vals[7,0] = 15.857855711422845
$vals[8,0] = 18.28061122244489
$vals[9,0] = 16.386456913827654
$vals[10,0] = 44.011623246492988
$vals[11,0] = 43.822581162324646
$vals[12,0] = 44.20066533066133
$vals[13,0] = 44.602379759519032
$vals[14,0] = 43.184564128256518
$vals[15,0] = 43.98799298597195
$vals[16,0] = 40.443453907815623
$vals[17,0] = 46.587321643286586
$vals[18,0] = 52.683928857715422
$vals[19,0] = 44.318816633266529
$vals[20,0] = 43.468127254509021
$vals[21,0] = 42.995522044088176
$vals[22,0] = 41.648597194388778
$vals[23,0] = 44.082514028056103
$vals[24,0] = 43.397236472945892
$vals[25,0] = 43.231824649298602
$vals[26,0] = 45.358548096192379
$vals[27,0] = 50.273642284569142
$vals[28,0] = 40.443453907815623
$vals[29,0] = 52.565765531062127
$vals[30,0] = 50.797472945891784
$vals[31,0] = 55.681328657314637
$vals[32,0] = 51.274631262525048
$vals[33,0] = 50.657132264529061
$vals[34,0] = 50.095769539078162
$vals[35,0] = 48.720430861723443
$vals[36,0] = 59.386322645290576
$vals[37,0] = 52.397356713426852
$vals[38,0] = 52.762242484969939
$vals[39,0] = 53.267468937875741
$vals[40,0] = 52.593833667334671
$vals[41,0] = 51.976334669338684
$vals[42,0] = 50.993949899799603
$vals[43,0] = 53.884967935871749
$vals[44,0] = 53.379741482965926
$vals[45,0] = 54.109513026052106
$vals[46,0] = 51.948266533066139
$vals[47,0] = 48.748498997995988
$vals[48,0] = 22.693649298597194
$vals[49,0] = 22.722418837675352
$vals[50,0] = 22.600148296593193
$vals[51,0] = 22.298068136272544
$vals[52,0] = 22.377184368737474
$vals[53,0] = 22.211759519038072
$vals[54,0] = 21.808985971943887
$vals[55,0] = 22.175797595190382
$vals[56,0] = 25.714450901803602
$vals[57,0] = 22.758380761523046
$vals[58,0] = 22.909420841683367
$vals[59,0] = 22.34122244488978
$vals[60,0] = 22.30526052104209
$vals[61,0] = 23.434464929859715
$vals[62,0] = 21.722677354709422
$ws.Range("G11:G73").Value = $vals

$vals = New-Object 'object[,]' 40,1
$vals[0,0] = 25.52744889779559
$vals[1,0] = 22.787150300601208
$vals[2,0] = 49.816354709418839
$vals[3,0] = 49.570202404809628
$vals[4,0] = 50.38250501002004
$vals[5,0] = 46.000993987975946
$vals[6,0] = 50.259428857715427
$vals[7,0] = 47.920981963927851
$vals[8,0] = 46.7394509018036
$vals[9,0] = 50.210198396793587
$vals[10,0] = 59.736292585170347
$vals[11,0] = 48.93020641282564
$vals[12,0] = 54.394787575150296
$vals[13,0] = 48.290210420841674
$vals[14,0] = 50.357889779559123
$vals[15,0] = 53.065565130260516
$vals[16,0] = 47.010218436873743
$vals[17,0] = 56.807080160320631
$vals[18,0] = 53.680945891783566
$vals[19,0] = 55.280935871743502
$vals[20,0] = 51.687112224448882
$vals[21,0] = 54.78254308617236
$vals[22,0] = 53.071440881763522
$vals[23,0] = 53.3346873747495
$vals[24,0] = 53.518959919839681
$vals[25,0] = 53.361012024048094
$vals[26,0] = 52.702895791583167
$vals[27,0] = 53.413661322645297
$vals[28,0] = 53.124090180360724
$vals[29,0] = 53.124090180360724
$vals[30,0] = 57.730903807615242
$vals[31,0] = 55.335360721442889
$vals[32,0] = 56.335697394789584
$vals[33,0] = 54.098102204408825
$vals[34,0] = 53.913829659318637
$vals[35,0] = 53.808531062124253
$vals[36,0] = 53.650583166332666
$vals[37,0] = 57.493981963927851
$vals[38,0] = 56.283048096192388
$vals[39,0] = 54.993140280561136
$ws.Range("G76:G115").Value = $vals

# Update the active cell/selection to match the saved workbook state.
$null = $ws.Range("I2").Select()
